# Update cryptos list (prices / 1h volume %) to the latest scrape,
# per "Updated cryptos list on Sat Sep 16 04:59:40 UTC 2023 with GitHub Actions".
#
# Notes:
#  - Column D ("Price") holds numeric-looking strings (e.g. "216.18"), but the
#    sheet stores them as TEXT (t="inlineStr"/shared string), not numbers -
#    some even use "." as a thousands separator (e.g. "26.712.00"), which is
#    not a valid number anyway. Assigning a clean single-decimal numeric
#    string (e.g. "216.18") to .Value would normally get auto-coerced to a
#    real number by Excel (and lose precision, e.g. "216.18000000000001").
#    To keep those cells as TEXT - exactly like the source workbook - a
#    leading apostrophe is used for those specific values, which is the
#    standard Excel "force text" input convention.
#  - Rows 13 and 14 swap coin identity (Polkadot <-> WrappedEther) in
#    addition to their price/volume values, so Name/Link/Price/Volume are
#    all rewritten for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '26.712.00'
$ws.Range('E2').Value = '  +0.04%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.646.87'
$ws.Range('E3').Value = '  +0.65%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.30%  '

# Row 5 - BNB
$ws.Range('D5').Value = "'216.18"
$ws.Range('E5').Value = '  +1.28%  '

# Row 6 - XRP
$ws.Range('D6').Value = "'0.504"
$ws.Range('E6').Value = '  -0.64%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  +0.25%  '

# Row 8 - Cardano
$ws.Range('E8').Value = '  -0.36%  '

# Row 9 - Dogecoin
$ws.Range('E9').Value = '  +0.57%  '

# Row 10 - Solana
$ws.Range('D10').Value = "'19.44"
$ws.Range('E10').Value = '  +0.99%  '

# Row 11 - TRON
$ws.Range('E11').Value = '  +0.22%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.878.18'
$ws.Range('E12').Value = '  +0.74%  '

# Row 13 - was Polkadot, now WrappedEther
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.693.83'
$ws.Range('E13').Value = '  +3.49%  '

# Row 14 - was WrappedEther, now Polkadot
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'4.24"
$ws.Range('E14').Value = '  +3.37%  '

# Row 15 - Polygon
$ws.Range('E15').Value = '  +1.45%  '

# Row 16 - Litecoin
$ws.Range('D16').Value = "'66.37"
$ws.Range('E16').Value = '  +4.81%  '

# Row 17 - WrappedBTC
$ws.Range('D17').Value = '26.771.64'
$ws.Range('E17').Value = '  +0.31%  '

# Row 18 - ShibaInu
$ws.Range('D18').Value = '0.0₃0756'
$ws.Range('E18').Value = '  +1.37%  '

# Row 19 - BitcoinCash
$ws.Range('D19').Value = "'220.42"
$ws.Range('E19').Value = '  +0.69%  '

# Row 20 - Dai
$ws.Range('E20').Value = '  +0.31%  '

# Row 21 - Uniswap
$ws.Range('D21').Value = "'4.40"
$ws.Range('E21').Value = '  +2.05%  '

# Row 22 - Chainlink
$ws.Range('E22').Value = '  +2.03%  '

# Row 23 - Avalanche
$ws.Range('E23').Value = '  +1.56%  '

# Row 24 - Toncoin
$ws.Range('E24').Value = '  +9.04%  '

# Row 25 - Monero
$ws.Range('D25').Value = "'147.07"
$ws.Range('E25').Value = '  -0.81%  '

# Row 26 - BinanceUSD
$ws.Range('E26').Value = '  +0.09%  '

# Row 27 - Stellar
$ws.Range('E27').Value = '  -0.84%  '

# Row 28 - Cosmos
$ws.Range('E28').Value = '  +2.63%  '

# Row 29 - EthereumClassic
$ws.Range('D29').Value = "'15.92"
$ws.Range('E29').Value = '  +2.64%  '

# Row 30 - Hedera
$ws.Range('E30').Value = '  +1.82%  '

# Row 31 - PancakeSwap
$ws.Range('E31').Value = '  +0.66%  '

# Row 32 - Filecoin
$ws.Range('E32').Value = '  +2.78%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range('E33').Value = '  +2.90%  '

# Row 34 - Maker
$ws.Range('D34').Value = '1.292.49'
$ws.Range('E34').Value = '  +7.85%  '

# Row 35 - LidoDAOToken
$ws.Range('E35').Value = '  +2.73%  '

# Row 36 - VeChain
$ws.Range('D36').Value = "'0.0185"
$ws.Range('E36').Value = '  +6.77%  '

# Row 37 - HuobiToken
$ws.Range('D37').Value = "'2.41"
$ws.Range('E37').Value = '  +0.76%  '

# Row 38 - ARBITRUM
$ws.Range('D38').Value = "'0.832"
$ws.Range('E38').Value = '  +2.49%  '

# Row 39 - ImmutableX
$ws.Range('D39').Value = "'0.528"
$ws.Range('E39').Value = '  +4.24%  '

# Row 40 - PaxDollar
$ws.Range('E40').Value = '  +0.24%  '

# Row 41 - TrustWalletToken
$ws.Range('D41').Value = "'0.811"
$ws.Range('E41').Value = '  +2.08%  '

# Row 42 - MXToken
$ws.Range('E42').Value = '  -1.86%  '

# Row 43 - FraxShare
$ws.Range('E43').Value = '  +0.43%  '

# Row 44 - RocketPoolETH
$ws.Range('D44').Value = '1.789.12'
$ws.Range('E44').Value = '  +0.94%  '

# Row 45 - Quant
$ws.Range('E45').Value = '  +1.45%  '

# Row 46 - Aave
$ws.Range('D46').Value = "'60.84"
$ws.Range('E46').Value = '  +10.84%  '

# Row 47 - RenderToken
$ws.Range('E47').Value = '  +3.58%  '

# Row 48 - Cronos
$ws.Range('E48').Value = '  +0.77%  '

# Row 49 - EnergySwap
$ws.Range('D49').Value = "'7.81"
$ws.Range('E49').Value = '  +2.29%  '

# Row 50 - Algorand
$ws.Range('D50').Value = "'0.0982"
$ws.Range('E50').Value = '  +3.44%  '

# Row 51 - Mantle
$ws.Range('E51').Value = '  -0.46%  '
